# Loan RBI, Variable Instalments
#
# Inserts a new (blank) column into the "Repayment Schedule" sheet just
# before the existing "Late" column (column N), pushing the old "Late"
# column to O and the old "Outstanding" column from P to Q. Also updates
# the active sheet/selection so that "Repayment Schedule" (not
# "Transactions") is the tab that is active/selected when the workbook is
# reopened.

$wb = $excel.ActiveWorkbook

$wsRepayment = $wb.Worksheets.Item("Repayment Schedule")
$wsTransactions = $wb.Worksheets.Item("Transactions")

# Insert a brand-new blank column at N; everything from N onward (N -> O,
# O -> P, P -> Q, ...) shifts one column to the right, matching an
# Excel "Insert" on the column header.
$wsRepayment.Columns("N").Insert()

# "Transactions" is no longer the active tab, so give it a plain
# (non-tab-selected) view with its own last selection.
$wsTransactions.Activate()
$wsTransactions.Range("B3").Select()

# Make "Repayment Schedule" the active sheet/tab, with R7 selected there.
# This must run last so it ends up as the workbook's active/selected tab.
$wsRepayment.Activate()
$wsRepayment.Range("R7").Select()
